# Adds a new "RecNumber"-style helper column (RecSource) to the
# "Not Normalized" sheet and renames the header cells so they match the
# field names used by the normalized import (CourseName / CourseDate /
# Completed), mirroring the changes made in NewFunctions.script and
# WorkbookCommandAndFunction.script.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Not Normalized")

# Rename headers to the "programmatic" field names.
$ws.Range("C3").Value = "CourseName"
$ws.Range("E3").Value = "CourseDate"
$ws.Range("F3").Value = "Completed"

# New column G: constant record-source marker (2) for every data row.
$lastRow = 11
for ($r = 4; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 7).Value = 2
}

# Leave the selection on the newly added column header, matching the
# workbook state saved by the author.
$ws.Range("G3").Select()

$wb.Save()
